$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: the phone number in A3 was stored as text; convert it to a real
# numeric value (matches the other rows, e.g. A2).
$ws.Range("A3").Value = 79174445

# Row 4: new redemption record - 79174445 redeemed 20 points.
# A4 keeps the legacy text representation of the phone number (leading
# apostrophe forces Excel to store it as text instead of auto-converting
# it to a number).
$ws.Range("A4").Value = "'79174445"
$ws.Range("B4").Value = 20
$ws.Range("C4").Value = "2025-08-18T08:51:52"
